$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 213.8864531954784
$ws.Range("C2").Value = 5729.65228008093
$ws.Range("D2").Value = 5909.953657128959
